$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Thbs1"
$ws.Range("C2").Value = "Itga4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 18.82499266666667
$ws.Range("H2").Value = 56.474978
$ws.Range("I2").Value = 0.06886869772378311
$ws.Range("J2").Value = 0.0688686977237831
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3809099999999999
$ws.Range("N2").Value = 1.14273
$ws.Range("O2").Value = 0.006635732896411959
$ws.Range("P2").Value = 0.006635732896411961
$ws.Range("Q2").Value = 7.170627956659999
$ws.Range("R2").Value = 64.53565160993999
$ws.Range("S2").Value = 0.0004569942830187589
$ws.Range("T2").Value = 0.000456994283018759

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Thbs1"
$ws.Range("C3").Value = "Itga4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 18.82499266666667
$ws.Range("H3").Value = 56.474978
$ws.Range("I3").Value = 0.06886869772378311
$ws.Range("J3").Value = 0.0688686977237831
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3194813333333333
$ws.Range("N3").Value = 0.9584440000000001
$ws.Range("O3").Value = 0.005565600255676026
$ws.Range("P3").Value = 0.005565600255676028
$ws.Range("Q3").Value = 6.014233757136889
$ws.Range("R3").Value = 54.12810381423201
$ws.Range("S3").Value = 0.0003832956416595622
$ws.Range("T3").Value = 0.0003832956416595623

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Thbs1"
$ws.Range("C4").Value = "Itga4"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 18.82499266666667
$ws.Range("H4").Value = 56.474978
$ws.Range("I4").Value = 0.06886869772378311
$ws.Range("J4").Value = 0.0688686977237831
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 54.69403966666666
$ws.Range("N4").Value = 164.082119
$ws.Range("O4").Value = 0.9528104755815301
$ws.Range("P4").Value = 0.9528104755815303
$ws.Range("Q4").Value = 1029.614895635376
$ws.Range("R4").Value = 9266.534060718381
$ws.Range("S4").Value = 0.06561881663087843
$ws.Range("T4").Value = 0.06561881663087843

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Thbs1"
$ws.Range("C5").Value = "Itga4"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 18.82499266666667
$ws.Range("H5").Value = 56.474978
$ws.Range("I5").Value = 0.06886869772378311
$ws.Range("J5").Value = 0.0688686977237831
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.008422
$ws.Range("N5").Value = 6.025265999999999
$ws.Range("O5").Value = 0.03498819126638183
$ws.Range("P5").Value = 0.03498819126638184
$ws.Range("Q5").Value = 37.808529421572
$ws.Range("R5").Value = 340.2767647941479
$ws.Range("S5").Value = 0.002409591168226358
$ws.Range("T5").Value = 0.002409591168226358

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Thbs1"
$ws.Range("C6").Value = "Itga4"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 121.8208923333333
$ws.Range("H6").Value = 365.462677
$ws.Range("I6").Value = 0.4456653109566078
$ws.Range("J6").Value = 0.4456653109566078
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.3809099999999999
$ws.Range("N6").Value = 1.14273
$ws.Range("O6").Value = 0.006635732896411959
$ws.Range("P6").Value = 0.006635732896411961
$ws.Range("Q6").Value = 46.40279609869
$ws.Range("R6").Value = 417.62516488821
$ws.Range("S6").Value = 0.002957315964704427
$ws.Range("T6").Value = 0.002957315964704428

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Thbs1"
$ws.Range("C7").Value = "Itga4"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 121.8208923333333
$ws.Range("H7").Value = 365.462677
$ws.Range("I7").Value = 0.4456653109566078
$ws.Range("J7").Value = 0.4456653109566078
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.3194813333333333
$ws.Range("N7").Value = 0.9584440000000001
$ws.Range("O7").Value = 0.005565600255676026
$ws.Range("P7").Value = 0.005565600255676028
$ws.Range("Q7").Value = 38.91950111050978
$ws.Range("R7").Value = 350.275509994588
$ws.Range("S7").Value = 0.002480394968606032
$ws.Range("T7").Value = 0.002480394968606033

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Thbs1"
$ws.Range("C8").Value = "Itga4"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 121.8208923333333
$ws.Range("H8").Value = 365.462677
$ws.Range("I8").Value = 0.4456653109566078
$ws.Range("J8").Value = 0.4456653109566078
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 54.69403966666666
$ws.Range("N8").Value = 164.082119
$ws.Range("O8").Value = 0.9528104755815301
$ws.Range("P8").Value = 0.9528104755815303
$ws.Range("Q8").Value = 6662.876717508063
$ws.Range("R8").Value = 59965.89045757256
$ws.Range("S8").Value = 0.424634576882756
$ws.Range("T8").Value = 0.424634576882756

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Thbs1"
$ws.Range("C9").Value = "Itga4"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 121.8208923333333
$ws.Range("H9").Value = 365.462677
$ws.Range("I9").Value = 0.4456653109566078
$ws.Range("J9").Value = 0.4456653109566078
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.008422
$ws.Range("N9").Value = 6.025265999999999
$ws.Range("O9").Value = 0.03498819126638183
$ws.Range("P9").Value = 0.03498819126638184
$ws.Range("Q9").Value = 244.667760221898
$ws.Range("R9").Value = 2202.009841997082
$ws.Range("S9").Value = 0.01559302314054133
$ws.Range("T9").Value = 0.01559302314054133

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Thbs1"
$ws.Range("C10").Value = "Itga4"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 87.673585
$ws.Range("H10").Value = 263.020755
$ws.Range("I10").Value = 0.3207419907481189
$ws.Range("J10").Value = 0.3207419907481188
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.3809099999999999
$ws.Range("N10").Value = 1.14273
$ws.Range("O10").Value = 0.006635732896411959
$ws.Range("P10").Value = 0.006635732896411961
$ws.Range("Q10").Value = 33.39574526234999
$ws.Range("R10").Value = 300.56170736115
$ws.Range("S10").Value = 0.002128358179267952
$ws.Range("T10").Value = 0.002128358179267953

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Thbs1"
$ws.Range("C11").Value = "Itga4"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 87.673585
$ws.Range("H11").Value = 263.020755
$ws.Range("I11").Value = 0.3207419907481189
$ws.Range("J11").Value = 0.3207419907481188
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.3194813333333333
$ws.Range("N11").Value = 0.9584440000000001
$ws.Range("O11").Value = 0.005565600255676026
$ws.Range("P11").Value = 0.005565600255676028
$ws.Range("Q11").Value = 28.01007383391334
$ws.Range("R11").Value = 252.09066450522
$ws.Range("S11").Value = 0.001785121705713768
$ws.Range("T11").Value = 0.001785121705713768

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Thbs1"
$ws.Range("C12").Value = "Itga4"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 87.673585
$ws.Range("H12").Value = 263.020755
$ws.Range("I12").Value = 0.3207419907481189
$ws.Range("J12").Value = 0.3207419907481188
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 54.69403966666666
$ws.Range("N12").Value = 164.082119
$ws.Range("O12").Value = 0.9528104755815301
$ws.Range("P12").Value = 0.9528104755815303
$ws.Range("Q12").Value = 4795.222535708871
$ws.Range("R12").Value = 43157.00282137984
$ws.Range("S12").Value = 0.3056063287436819
$ws.Range("T12").Value = 0.3056063287436819

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Thbs1"
$ws.Range("C13").Value = "Itga4"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 87.673585
$ws.Range("H13").Value = 263.020755
$ws.Range("I13").Value = 0.3207419907481189
$ws.Range("J13").Value = 0.3207419907481188
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.008422
$ws.Range("N13").Value = 6.025265999999999
$ws.Range("O13").Value = 0.03498819126638183
$ws.Range("P13").Value = 0.03498819126638184
$ws.Range("Q13").Value = 176.08555693287
$ws.Range("R13").Value = 1584.77001239583
$ws.Range("S13").Value = 0.01122218211945525
$ws.Range("T13").Value = 0.01122218211945525

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Thbs1"
$ws.Range("C14").Value = "Itga4"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 45.02666966666666
$ws.Range("H14").Value = 135.080009
$ws.Range("I14").Value = 0.1647240005714903
$ws.Range("J14").Value = 0.1647240005714903
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.3809099999999999
$ws.Range("N14").Value = 1.14273
$ws.Range("O14").Value = 0.006635732896411959
$ws.Range("P14").Value = 0.006635732896411961
$ws.Range("Q14").Value = 17.15110874273
$ws.Range("R14").Value = 154.35997868457
$ws.Range("S14").Value = 0.00109306446942082
$ws.Range("T14").Value = 0.001093064469420821

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Thbs1"
$ws.Range("C15").Value = "Itga4"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 45.02666966666666
$ws.Range("H15").Value = 135.080009
$ws.Range("I15").Value = 0.1647240005714903
$ws.Range("J15").Value = 0.1647240005714903
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.3194813333333333
$ws.Range("N15").Value = 0.9584440000000001
$ws.Range("O15").Value = 0.005565600255676026
$ws.Range("P15").Value = 0.005565600255676028
$ws.Range("Q15").Value = 14.38518046066622
$ws.Range("R15").Value = 129.466624145996
$ws.Range("S15").Value = 0.0009167879396966642
$ws.Range("T15").Value = 0.0009167879396966643

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Thbs1"
$ws.Range("C16").Value = "Itga4"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 45.02666966666666
$ws.Range("H16").Value = 135.080009
$ws.Range("I16").Value = 0.1647240005714903
$ws.Range("J16").Value = 0.1647240005714903
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 54.69403966666666
$ws.Range("N16").Value = 164.082119
$ws.Range("O16").Value = 0.9528104755815301
$ws.Range("P16").Value = 0.9528104755815303
$ws.Range("Q16").Value = 2462.690456806563
$ws.Range("R16").Value = 22164.21411125907
$ws.Range("S16").Value = 0.1569507533242139
$ws.Range("T16").Value = 0.1569507533242139

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Thbs1"
$ws.Range("C17").Value = "Itga4"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 45.02666966666666
$ws.Range("H17").Value = 135.080009
$ws.Range("I17").Value = 0.1647240005714903
$ws.Range("J17").Value = 0.1647240005714903
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 2.008422
$ws.Range("N17").Value = 6.025265999999999
$ws.Range("O17").Value = 0.03498819126638183
$ws.Range("P17").Value = 0.03498819126638184
$ws.Range("Q17").Value = 90.43255394526599
$ws.Range("R17").Value = 813.8929855073939
$ws.Range("S17").Value = 0.005763394838158892
$ws.Range("T17").Value = 0.005763394838158892

Write-Host "Edit complete"